# Actualización desde MV -datos-
# Add a new column (BH) with header "Agosto.2021" that carries forward the
# last known value of each row (same as column BG), mirroring a new
# reporting period being appended to the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column BH (60), row 1 -- same formatting as BG1
$ws.Range("BH1").Value = "Agosto.2021"
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122) # xlPasteFormats

# Carry forward the last value of each data row (2-19) from column BG into BH
for ($r = 2; $r -le 19; $r++) {
    $lastValue = $ws.Cells.Item($r, 59).Value2
    $ws.Cells.Item($r, 60).Value = $lastValue
}
